$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.208.69'
$ws.Range("E2").Value = '  -2.44%  '
$ws.Range("D3").Value = '2.870.42'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.75'
$ws.Range("E5").Value = '  -4.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.64'
$ws.Range("E6").Value = '  -3.18%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").Value = '2.866.86'
$ws.Range("E9").Value = '  -2.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.84'
$ws.Range("E10").Value = '  -6.26%  '
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.429'
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.72'
$ws.Range("E14").Value = '  -3.33%  '
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").Value = '3.349.03'
$ws.Range("E16").Value = '  -2.40%  '
$ws.Range("D17").Value = '61.259.43'
$ws.Range("E17").Value = '  -2.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.52'
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").Value = '2.866.65'
$ws.Range("E19").Value = '  -2.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '430.08'
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.03'
$ws.Range("E21").Value = '  -2.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.650'
$ws.Range("E22").Value = '  -2.40%  '
$ws.Range("E23").Value = '  -3.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.82'
$ws.Range("E24").Value = '  -3.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.60'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.89'
$ws.Range("E27").Value = '  -10.92%  '
$ws.Range("E28").Value = '  -6.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000103'
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.96'
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.48'
$ws.Range("E31").Value = '  -4.77%  '
$ws.Range("E32").Value = '  -8.88%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  -3.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.41'
$ws.Range("E35").Value = '  -3.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.953'
$ws.Range("E37").Value = '  -4.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.80'
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  -5.24%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.81'
$ws.Range("E40").Value = '  -10.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.19'
$ws.Range("E41").Value = '  -3.42%  '
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.15'
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.266'
$ws.Range("E44").Value = '  -5.16%  '
$ws.Range("D45").Value = '2.677.15'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.70'
$ws.Range("E46").Value = '  -1.87%  '
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '336.64'
$ws.Range("E49").Value = '  -7.13%  '
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.34'
$ws.Range("E51").Value = '  -6.57%  '
